$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 394, pushing the existing rows 394-495
# down to 396-497 (dimension grows from A1:R495 to A1:R497).
$ws.Rows.Item(394).Resize(2).Insert()

# Row 394 (new) - "Primera" quality record for 2021-12-21 (serial 44551)
$ws.Range("A394").Value = 8
$ws.Range("B394").Value = "Terminal La Palmera de La Serena"
$ws.Range("C394").Value = "Coquimbo"
$ws.Range("D394").Value = 44551
$ws.Range("E394").Value = 4
$ws.Range("F394").Value = 100112008
$ws.Range("G394").Value = "Coliflor"
$ws.Range("H394").Value = "Sin especificar"
$ws.Range("I394").Value = "Primera"
$ws.Range("J394").Value = 2300
$ws.Range("K394").Value = 600
$ws.Range("L394").Value = 700
$ws.Range("M394").Value = 650
$ws.Range("N394").Value = "`$/unidad"
$ws.Range("O394").Value = "Provincia del Elquí"
$ws.Range("P394").Value = 650
$ws.Range("Q394").Value = 1
$ws.Range("R394").Value = "Hortaliza"

# Row 395 (new) - "Segunda" quality record for 2021-12-21 (serial 44551)
$ws.Range("A395").Value = 8
$ws.Range("B395").Value = "Terminal La Palmera de La Serena"
$ws.Range("C395").Value = "Coquimbo"
$ws.Range("D395").Value = 44551
$ws.Range("E395").Value = 4
$ws.Range("F395").Value = 100112008
$ws.Range("G395").Value = "Coliflor"
$ws.Range("H395").Value = "Sin especificar"
$ws.Range("I395").Value = "Segunda"
$ws.Range("J395").Value = 1200
$ws.Range("K395").Value = 500
$ws.Range("L395").Value = 550
$ws.Range("M395").Value = 525
$ws.Range("N395").Value = "`$/unidad"
$ws.Range("O395").Value = "Provincia del Elquí"
$ws.Range("P395").Value = 525
$ws.Range("Q395").Value = 1
$ws.Range("R395").Value = "Hortaliza"
